$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @(2, 7, 1.699314666666667),
  @(2, 8, 5.097944),
  @(2, 9, 0.01313494849782423),
  @(2, 10, 0.0136763907731691),
  @(2, 13, 2.906846333333333),
  @(2, 14, 8.720538999999999),
  @(2, 15, 0.005520525738044089),
  @(2, 16, 0.005624540846623205),
  @(2, 17, 4.939646607979554),
  @(2, 18, 44.45681947181599),
  @(2, 19, 0.00007251182125012223),
  @(2, 20, 0.00007692341853807031),
  @(3, 7, 1.699314666666667),
  @(3, 8, 5.097944),
  @(3, 9, 0.01313494849782423),
  @(3, 10, 0.0136763907731691),
  @(3, 15, 0.3528665483720876),
  @(3, 16, 0.3595150912979765),
  @(3, 17, 315.7373285525502),
  @(3, 18, 2841.635956972952),
  @(3, 19, 0.004634883939472375),
  @(3, 20, 0.004916868877442692),
  @(4, 7, 1.699314666666667),
  @(4, 8, 5.097944),
  @(4, 9, 0.01313494849782423),
  @(4, 10, 0.0136763907731691),
  @(4, 13, 137.0717086666666),
  @(4, 14, 411.2151259999999),
  @(4, 15, 0.2603191943704447),
  @(4, 16, 0.2652240042658267),
  @(4, 17, 232.9279649223271),
  @(4, 18, 2096.351684300944),
  @(4, 19, 0.003419279211050887),
  @(4, 20, 0.003627307124764113),
  @(5, 7, 1.699314666666667),
  @(5, 8, 5.097944),
  @(5, 9, 0.01313494849782423),
  @(5, 10, 0.0136763907731691),
  @(5, 13, 29.2127365),
  @(5, 14, 58.425473),
  @(5, 15, 0.05547925319534149),
  @(5, 16, 0.03768304451958546),
  @(5, 17, 49.64163158791867),
  @(5, 18, 297.849789527512),
  @(5, 19, 0.000728717133418561),
  @(5, 20, 0.0005153680423725789),
  @(6, 7, 1.699314666666667),
  @(6, 8, 5.097944),
  @(6, 9, 0.01313494849782423),
  @(6, 10, 0.0136763907731691),
  @(6, 13, 171.5584106666666),
  @(6, 14, 514.6752319999999),
  @(6, 15, 0.3258144783240821),
  @(6, 16, 0.331953319069988),
  @(6, 17, 291.5317234358898),
  @(6, 18, 2623.785510923008),
  @(6, 19, 0.004279556392632289),
  @(6, 20, 0.004539923310051642),
  @(7, 9, 0.7574068660664961),
  @(7, 10, 0.7886283129562799),
  @(7, 13, 2.906846333333333),
  @(7, 14, 8.720538999999999),
  @(7, 15, 0.005520525738044089),
  @(7, 16, 0.005624540846623205),
  @(7, 17, 284.8372231870975),
  @(7, 18, 2563.535008683878),
  @(7, 19, 0.004181284098291403),
  @(7, 20, 0.004435672159026145),
  @(8, 9, 0.7574068660664961),
  @(8, 10, 0.7886283129562799),
  @(8, 15, 0.3528665483720876),
  @(8, 16, 0.3595150912979765),
  @(8, 19, 0.2672635465422046),
  @(8, 20, 0.2835237799326462),
  @(9, 9, 0.7574068660664961),
  @(9, 10, 0.7886283129562799),
  @(9, 13, 137.0717086666666),
  @(9, 14, 411.2151259999999),
  @(9, 15, 0.2603191943704447),
  @(9, 16, 0.2652240042658267),
  @(9, 17, 13431.43750889394),
  @(9, 18, 120882.9375800454),
  @(9, 19, 0.1971675451850735),
  @(9, 20, 0.2091631590396681),
  @(10, 9, 0.7574068660664961),
  @(10, 10, 0.7886283129562799),
  @(10, 13, 29.2127365),
  @(10, 14, 58.425473),
  @(10, 15, 0.05547925319534149),
  @(10, 16, 0.03768304451958546),
  @(10, 17, 2862.509328731757),
  @(10, 18, 17175.05597239054),
  @(10, 19, 0.04202036729439323),
  @(10, 20, 0.02971791582653707),
  @(11, 9, 0.7574068660664961),
  @(11, 10, 0.7886283129562799),
  @(11, 13, 171.5584106666666),
  @(11, 14, 514.6752319999999),
  @(11, 15, 0.3258144783240821),
  @(11, 16, 0.331953319069988),
  @(11, 17, 16810.73428214187),
  @(11, 18, 151296.6085392768),
  @(11, 19, 0.2467741229465334),
  @(11, 20, 0.2617877859984024),
  @(12, 7, 3.143799666666666),
  @(12, 8, 9.431398999999999),
  @(12, 9, 0.02430017672368134),
  @(12, 10, 0.02530186645080374),
  @(12, 13, 2.906846333333333),
  @(12, 14, 8.720538999999999),
  @(12, 15, 0.005520525738044089),
  @(12, 16, 0.005624540846623205),
  @(12, 17, 9.138542533784554),
  @(12, 18, 82.24688280406097),
  @(12, 19, 0.0001341497510421027),
  @(12, 20, 0.000142311381348351),
  @(13, 7, 3.143799666666666),
  @(13, 8, 9.431398999999999),
  @(13, 9, 0.02430017672368134),
  @(13, 10, 0.02530186645080374),
  @(13, 15, 0.3528665483720876),
  @(13, 16, 0.3595150912979765),
  @(13, 17, 584.1266057008851),
  @(13, 18, 5257.139451307967),
  @(13, 19, 0.00857471948531718),
  @(13, 20, 0.009096402827069918),
  @(14, 7, 3.143799666666666),
  @(14, 8, 9.431398999999999),
  @(14, 9, 0.02430017672368134),
  @(14, 10, 0.02530186645080374),
  @(14, 13, 137.0717086666666),
  @(14, 14, 411.2151259999999),
  @(14, 15, 0.2603191943704447),
  @(14, 16, 0.2652240042658267),
  @(14, 17, 430.925992015697),
  @(14, 18, 3878.333928141273),
  @(14, 19, 0.006325802427768159),
  @(14, 20, 0.00671066233548135),
  @(15, 7, 3.143799666666666),
  @(15, 8, 9.431398999999999),
  @(15, 9, 0.02430017672368134),
  @(15, 10, 0.02530186645080374),
  @(15, 13, 29.2127365),
  @(15, 14, 58.425473),
  @(15, 15, 0.05547925319534149),
  @(15, 16, 0.03768304451958546),
  @(15, 17, 91.83899127112115),
  @(15, 18, 551.0339476267269),
  @(15, 19, 0.001348155657144661),
  @(15, 20, 0.0009534513598942432),
  @(16, 7, 3.143799666666666),
  @(16, 8, 9.431398999999999),
  @(16, 9, 0.02430017672368134),
  @(16, 10, 0.02530186645080374),
  @(16, 13, 171.5584106666666),
  @(16, 14, 514.6752319999999),
  @(16, 15, 0.3258144783240821),
  @(16, 16, 0.331953319069988),
  @(16, 17, 539.3452742677297),
  @(16, 18, 4854.107468409567),
  @(16, 19, 0.00791734940240924),
  @(16, 20, 0.008399038547009881),
  @(17, 7, 15.365523),
  @(17, 8, 30.731046),
  @(17, 9, 0.1187686761057793),
  @(17, 10, 0.08244299936684966),
  @(17, 13, 2.906846333333333),
  @(17, 14, 8.720538999999999),
  @(17, 15, 0.005520525738044089),
  @(17, 16, 0.005624540846623205),
  @(17, 17, 44.66521419229899),
  @(17, 18, 267.991285153794),
  @(17, 19, 0.0006556655333153767),
  @(17, 20, 0.000463704017456977),
  @(18, 7, 15.365523),
  @(18, 8, 30.731046),
  @(18, 9, 0.1187686761057793),
  @(18, 10, 0.08244299936684966),
  @(18, 15, 0.3528665483720876),
  @(18, 16, 0.3595150912979765),
  @(18, 17, 2854.956341516953),
  @(18, 18, 17129.73804910172),
  @(18, 19, 0.0419094927921688),
  @(18, 20, 0.02963950244425198),
  @(19, 7, 15.365523),
  @(19, 8, 30.731046),
  @(19, 9, 0.1187686761057793),
  @(19, 10, 0.08244299936684966),
  @(19, 13, 137.0717086666666),
  @(19, 14, 411.2151259999999),
  @(19, 15, 0.2603191943704447),
  @(19, 16, 0.2652240042658267),
  @(19, 17, 2106.178492166965),
  @(19, 18, 12637.07095300179),
  @(19, 19, 0.03091776608030075),
  @(19, 20, 0.02186586241576088),
  @(20, 7, 15.365523),
  @(20, 8, 30.731046),
  @(20, 9, 0.1187686761057793),
  @(20, 10, 0.08244299936684966),
  @(20, 13, 29.2127365),
  @(20, 14, 58.425473),
  @(20, 15, 0.05547925319534149),
  @(20, 16, 0.03768304451958546),
  @(20, 17, 448.8689745836895),
  @(20, 18, 1795.475898334758),
  @(20, 19, 0.006589197453348036),
  @(20, 20, 0.003106703215469151),
  @(21, 7, 15.365523),
  @(21, 8, 30.731046),
  @(21, 9, 0.1187686761057793),
  @(21, 10, 0.08244299936684966),
  @(21, 13, 171.5584106666666),
  @(21, 14, 514.6752319999999),
  @(21, 15, 0.3258144783240821),
  @(21, 16, 0.331953319069988),
  @(21, 17, 2636.084704942112),
  @(21, 18, 15816.50822965267),
  @(21, 19, 0.03869655424664637),
  @(21, 20, 0.02736722727391067),
  @(22, 7, 11.176493),
  @(22, 8, 33.529479),
  @(22, 9, 0.08638933260621913),
  @(22, 10, 0.08995043045289769),
  @(22, 13, 2.906846333333333),
  @(22, 14, 8.720538999999999),
  @(22, 15, 0.005520525738044089),
  @(22, 16, 0.005624540846623205),
  @(22, 17, 32.48834769657567),
  @(22, 18, 292.395129269181),
  @(22, 19, 0.0004769145341450841),
  @(22, 20, 0.000505929870253663),
  @(23, 7, 11.176493),
  @(23, 8, 33.529479),
  @(23, 9, 0.08638933260621913),
  @(23, 10, 0.08995043045289769),
  @(23, 15, 0.3528665483720876),
  @(23, 16, 0.3595150912979765),
  @(23, 17, 2076.62307142229),
  @(23, 18, 18689.60764280061),
  @(23, 19, 0.03048390561292479),
  @(23, 20, 0.0323385372165658),
  @(24, 7, 11.176493),
  @(24, 8, 33.529479),
  @(24, 9, 0.08638933260621913),
  @(24, 10, 0.08995043045289769),
  @(24, 13, 137.0717086666666),
  @(24, 14, 411.2151259999999),
  @(24, 15, 0.2603191943704447),
  @(24, 16, 0.2652240042658267),
  @(24, 17, 1531.980992411039),
  @(24, 18, 13787.82893169935),
  @(24, 19, 0.02248880146625135),
  @(24, 20, 0.02385701335015229),
  @(25, 7, 11.176493),
  @(25, 8, 33.529479),
  @(25, 9, 0.08638933260621913),
  @(25, 10, 0.08995043045289769),
  @(25, 13, 29.2127365),
  @(25, 14, 58.425473),
  @(25, 15, 0.05547925319534149),
  @(25, 16, 0.03768304451958546),
  @(25, 17, 326.4959450030945),
  @(25, 18, 1958.975670018567),
  @(25, 19, 0.004792815657037001),
  @(25, 20, 0.003389606075312419),
  @(26, 7, 11.176493),
  @(26, 8, 33.529479),
  @(26, 9, 0.08638933260621913),
  @(26, 10, 0.08995043045289769),
  @(26, 13, 171.5584106666666),
  @(26, 14, 514.6752319999999),
  @(26, 15, 0.3258144783240821),
  @(26, 16, 0.331953319069988),
  @(26, 17, 1917.421375907125),
  @(26, 18, 17256.79238316413),
  @(26, 19, 0.0281468953358609),
  @(26, 20, 0.02985934394061351)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
